$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.252.11'
$ws.Range('E2').Value = '  +1.66%  '

$ws.Range('D3').Value = '1.657.65'
$ws.Range('E3').Value = '  +0.83%  '

$ws.Range('E4').Value = '  -0.67%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '220.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.62%  '

$ws.Range('E6').Value = '  +0.18%  '

$ws.Range('E7').Value = '  -0.70%  '

$ws.Range('E8').Value = '  +0.90%  '

$ws.Range('E9').Value = '  -0.30%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.63'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.26%  '

$ws.Range('E11').Value = '  +0.41%  '

$ws.Range('D12').Value = '1.886.54'

$ws.Range('D13').Value = '1.653.48'
$ws.Range('E13').Value = '  +0.53%  '

$ws.Range('E14').Value = '  +1.22%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.85'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.76%  '

$ws.Range('D17').Value = '27.211.46'
$ws.Range('E17').Value = '  +1.45%  '

$ws.Range('D18').Value = '0.0₃0739'
$ws.Range('E18').Value = '  +0.44%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '220.96'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.91%  '

$ws.Range('E20').Value = '  -0.61%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.74'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.67%  '

$ws.Range('E22').Value = '  +0.84%  '

$ws.Range('E23').Value = '  -0.61%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.29'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.67%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.79'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.93%  '

$ws.Range('E26').Value = '  -0.57%  '

$ws.Range('E27').Value = '  +2.73%  '

$ws.Range('E28').Value = '  +0.55%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.01'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.19%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0515'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.28%  '

$ws.Range('E31').Value = '  +1.60%  '

$ws.Range('E32').Value = '  +0.59%  '

$ws.Range('E33').Value = '  -0.13%  '

$ws.Range('E34').Value = '  +2.38%  '

$ws.Range('D35').Value = '1.273.08'
$ws.Range('E35').Value = '  -1.15%  '

$ws.Range('E36').Value = '  +0.07%  '

$ws.Range('E37').Value = '  -1.37%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.540'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.86%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.829'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.07%  '

$ws.Range('E40').Value = '  -0.57%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.809'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.19%  '

$ws.Range('E42').Value = '  +1.10%  '

$ws.Range('D43').Value = '1.799.15'
$ws.Range('E43').Value = '  +0.82%  '

$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.12'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.14%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.96'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.95%  '

$ws.Range('E46').Value = '  +0.80%  '

$ws.Range('E47').Value = '  +0.48%  '

$ws.Range('E48').Value = '  -0.60%  '

$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0977'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.81%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.65'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.12%  '

$ws.Range('E51').Value = '  +0.00%  '
